# CVX.xlsx update — "Add files via upload"
#
# Semantic changes:
#   1. Sheet1!D2 (share Price) bumped from 132.5 to 138. Everything else on
#      Sheet1 (D4, D7) and Sheet2 (Q23, Q24, Q25, etc.) cascades from this.
#   2. Sheet2 FCF-margin assumptions in row 28 (E28:N28) are switched from a
#      "+2%/yr" compounding formula off of D28 to flat, hard-keyed 10%
#      assumptions for every forecast year. This also ripples into the FCF
#      build in row 32 (E32:N32 and onward through the terminal-growth
#      block), the FCF/Capex multiple in row 33, and the cumulative net-cash
#      roll-forward in row 34.
#   3. Sheet2!D27 picks up the same percent number format used across row 28
#      (an incidental formatting artifact of whatever the author selected
#      when re-keying row 28), even though it stays empty.
#   4. Cosmetic: selections on both sheets move, and Sheet2's frozen pane
#      view scrolls down a bit. Sheet2 remains the active/visible sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- 1. Sheet1: bump the share price input -------------------------------
$ws1.Range("D2").Value = 138

# --- 2. Sheet2: flatten the FCF-margin assumption row to 10% -------------
$ws2.Range("E28:N28").Value = 0.1

# --- 3. Sheet2: pick up the percent format on the now-blank D27 ----------
$ws2.Range("D27").NumberFormat = "0%"

# --- 4. Cosmetic view/selection state -------------------------------------
# Sheet1's remembered selection moves to D3 (Sheet1 is not the active tab).
$ws1.Select()
$ws1.Range("D3").Select()

# Sheet2 stays the visible/active sheet; its bottom-right pane selection
# moves from N28 to F28, and the frozen view scrolls so row 5 leads.
$ws2.Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 2
$ws2.Range("F28").Select()

Write-Host "Applied CVX price/FCF-margin update."
